$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
  $r = $ws.Range($addr)
  $r.NumberFormat = "@"
  $r.Value = $val
  $r.Style = "Normal"
}

Set-TextValue "D2" "42.065.06"
Set-TextValue "E2" "  -1.47%  "
Set-TextValue "D3" "2.304.68"
Set-TextValue "E3" "  -1.82%  "
Set-TextValue "E4" "  -0.02%  "
Set-TextValue "D5" "315.50"
Set-TextValue "E5" "  -1.53%  "
Set-TextValue "D6" "104.46"
Set-TextValue "E6" "  -0.86%  "
Set-TextValue "E7" "  -1.33%  "
Set-TextValue "E8" "  -0.02%  "
Set-TextValue "D9" "0.608"
Set-TextValue "E9" "  -1.17%  "
Set-TextValue "D10" "39.72"
Set-TextValue "E10" "  -3.48%  "
Set-TextValue "D11" "0.0906"
Set-TextValue "E11" "  -1.71%  "
Set-TextValue "D12" "8.47"
Set-TextValue "E12" "  +0.93%  "
Set-TextValue "E13" "  +1.08%  "
Set-TextValue "D14" "0.974"
Set-TextValue "E14" "  -2.41%  "
Set-TextValue "D15" "15.43"
Set-TextValue "E15" "  -2.95%  "
Set-TextValue "D16" "2.653.42"
Set-TextValue "E16" "  -1.87%  "
Set-TextValue "D17" "2.310.74"
Set-TextValue "E17" "  -1.10%  "
Set-TextValue "D18" "42.064.65"
Set-TextValue "E18" "  -1.50%  "
Set-TextValue "D19" "7.69"
Set-TextValue "E19" "  -0.13%  "
Set-TextValue "E20" "  -0.26%  "
Set-TextValue "B21" "BitcoinCash"
Set-TextValue "C21" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D21" "286.60"
Set-TextValue "E21" "  +11.14%  "
Set-TextValue "B22" "Litecoin"
Set-TextValue "C22" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D22" "73.73"
Set-TextValue "E22" "  -4.46%  "
Set-TextValue "D23" "3.53"
Set-TextValue "E23" "  -3.09%  "
Set-TextValue "D24" "2.29"
Set-TextValue "E24" "  -0.32%  "
Set-TextValue "D25" "9.97"
Set-TextValue "E25" "  +6.10%  "
Set-TextValue "E26" "  +0.52%  "
Set-TextValue "D27" "10.94"
Set-TextValue "E27" "  -3.60%  "
Set-TextValue "D28" "23.53"
Set-TextValue "E28" "  +2.68%  "
Set-TextValue "E29" "  +2.30%  "
Set-TextValue "D30" "35.89"
Set-TextValue "D31" "165.56"
Set-TextValue "E31" "  -5.27%  "
Set-TextValue "D32" "0.0885"
Set-TextValue "E32" "  -0.15%  "
Set-TextValue "E33" "  -1.01%  "
Set-TextValue "D34" "5.89"
Set-TextValue "E34" "  -3.14%  "
Set-TextValue "B35" "Stellar"
Set-TextValue "C35" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D35" "0.132"
Set-TextValue "E35" "  +1.42%  "
Set-TextValue "B36" "Kaspa"
Set-TextValue "C36" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D36" "0.119"
Set-TextValue "E36" "  -1.97%  "
Set-TextValue "D37" "4.65"
Set-TextValue "E37" "  +1.34%  "
Set-TextValue "E38" "  +9.33%  "
Set-TextValue "D39" "0.0353"
Set-TextValue "E39" "  -1.72%  "
Set-TextValue "D40" "3.64"
Set-TextValue "E40" "  -3.14%  "
Set-TextValue "D41" "102.15"
Set-TextValue "E41" "  +18.91%  "
Set-TextValue "E42" "  +2.02%  "
Set-TextValue "D43" "70.76"
Set-TextValue "E43" "  -0.98%  "
Set-TextValue "D44" "0.226"
Set-TextValue "E44" "  -3.33%  "
Set-TextValue "E45" "  +0.22%  "
Set-TextValue "D46" "116.56"
Set-TextValue "E46" "  +1.52%  "
Set-TextValue "D47" "12.10"
Set-TextValue "E47" "  +1.75%  "
Set-TextValue "D48" "78.45"
Set-TextValue "E48" "  +7.20%  "
Set-TextValue "D49" "9.11"
Set-TextValue "E49" "  -0.03%  "
Set-TextValue "D50" "5.33"
Set-TextValue "D51" "1.28"
Set-TextValue "E51" "  +2.08%  "
